# Update in CLick event13
# Insert a new "SCROLL_DOWN" step row right before the existing row 14 (WAIT),
# which pushes the subsequent rows (WAIT/CLICK/... steps) down by one row.
# Also set the CSS ObjectType on the "CLICK MyaccountSection" row, and update
# the sheet view (topLeftCell / selection) to reflect scrolled position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC23_Verify_UserRegistration")

# Insert a new blank row above row 14 - existing rows 14.. shift down to 15..
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the SCROLL_DOWN keyword (column B),
# matching the thin-bordered look of every other step row on this sheet.
$ws.Cells.Item(14, 2).Value = "SCROLL_DOWN"
$newRowRange = $ws.Range("A14:E14")
$newRowRange.Borders.LineStyle = 1
$newRowRange.Borders.Weight = 2

# The row that previously held "CLICK / MyaccountSection" (now shifted to row 22)
# gets a CSS ObjectType value added in column D.
$ws.Cells.Item(22, 4).Value = "CSS"

# Update the sheet view to match the scrolled state captured in the workbook.
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("D16").Select()
